$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (I1, J1) - copy formatting from the existing header
# style (H1, xfId "1": bold/bordered/centered) so the new columns match
# the rest of row 1, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data row values (I2, J2), unstyled like the other data cells.
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9
